# Weekly update: insert 3 new rows of data (new week) above the existing
# block of rows 797-818, shifting the old rows down to 800-821, then fill
# in the new rows with the latest week's figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 797 onward down by 3 to make room for the new week's rows.
$ws.Rows("797:799").Insert()

# New row 797
$ws.Range("A797").Value = 12
$ws.Range("B797").Value = "Mapocho Venta Directa de Santiago"
$ws.Range("C797").Value = "Metropolitana"
$ws.Range("D797").Value = 44448
$ws.Range("E797").Value = 13
$ws.Range("F797").Value = 100112020
$ws.Range("G797").Value = "Tomate"
$ws.Range("H797").Value = "Larga vida"
$ws.Range("I797").Value = "Primera"
$ws.Range("J797").Value = 580
$ws.Range("K797").Value = 14000
$ws.Range("L797").Value = 15000
$ws.Range("M797").Value = 14517
$ws.Range("N797").Value = "$/bandeja 18 kilos"
$ws.Range("O797").Value = "Provincia de Quillota"
$ws.Range("P797").Value = 806
$ws.Range("Q797").Value = 18
$ws.Range("R797").Value = "Hortaliza"

# New row 798
$ws.Range("A798").Value = 12
$ws.Range("B798").Value = "Mapocho Venta Directa de Santiago"
$ws.Range("C798").Value = "Metropolitana"
$ws.Range("D798").Value = 44448
$ws.Range("E798").Value = 13
$ws.Range("F798").Value = 100112020
$ws.Range("G798").Value = "Tomate"
$ws.Range("H798").Value = "Larga vida"
$ws.Range("I798").Value = "Segunda"
$ws.Range("J798").Value = 420
$ws.Range("K798").Value = 10000
$ws.Range("L798").Value = 11000
$ws.Range("M798").Value = 10571
$ws.Range("N798").Value = "$/bandeja 18 kilos"
$ws.Range("O798").Value = "Provincia de Quillota"
$ws.Range("P798").Value = 587
$ws.Range("Q798").Value = 18
$ws.Range("R798").Value = "Hortaliza"

# New row 799
$ws.Range("A799").Value = 12
$ws.Range("B799").Value = "Mapocho Venta Directa de Santiago"
$ws.Range("C799").Value = "Metropolitana"
$ws.Range("D799").Value = 44448
$ws.Range("E799").Value = 13
$ws.Range("F799").Value = 100112020
$ws.Range("G799").Value = "Tomate"
$ws.Range("H799").Value = "Larga vida"
$ws.Range("I799").Value = "Segunda"
$ws.Range("J799").Value = 580
$ws.Range("K799").Value = 19000
$ws.Range("L799").Value = 20000
$ws.Range("M799").Value = 19517
$ws.Range("N799").Value = "$/bandeja 18 kilos"
$ws.Range("O799").Value = "Región de Arica y Parinacota"
$ws.Range("P799").Value = 1084
$ws.Range("Q799").Value = 18
$ws.Range("R799").Value = "Hortaliza"
